$wb = $excel.ActiveWorkbook

# --- Thermotaxis Index sheet: add thermotaxis gradient-quantification columns ---
$ws2 = $wb.Worksheets.Item("Thermotaxis Index ")

# Reorganize the header row: columns C:D keep their style, and the new
# gradient-quantification columns (C:H) get the same centered style as A:B.
# Existing "Tstart Camera" (row 7) is dropped - its info is now redundant.
$ws2.Range("C1").Value = "Low gradient"
$ws2.Range("D1").Value = "High gradient"
$ws2.Range("E1").Value = "Distance between"
$ws2.Range("F1").Value = "Gradient slope (value per per cm)"
$ws2.Range("G1").Value = "T(start)"
$ws2.Range("H1").Value = "pixelspercm"
$ws2.Range("I1").Value = "OdorXCoord"
$ws2.Range("J1").Value = "OdorYCoord"

# Match the centered header style already used by A1:D1 across the newly
# populated header cells E1:H1 (I1/J1 stay unstyled, like before).
$ws2.Range("E1:H1").HorizontalAlignment = -4108

# Drop the now-obsolete "Tstart Camera" row.
$ws2.Range("A7").ClearContents()

# Column F (Gradient slope...) needs to be wide enough to show the label.
$ws2.Columns.Item(6).ColumnWidth = 22.25

# Make this the active/selected sheet & cell (must be the last sheet
# activated so it ends up as the workbook's active tab).
$ws2.Activate()
$ws2.Range("F4").Select()
